# Generate Report for Handback
#
# - Flip the "Ready for handoff" status to "Handed back: in sync with en-US"
#   (Overview sheet + each language sheet).
# - Record the handback: for a.md/b.md rows in each language sheet, fill in
#   "Latest Target File" (E) and "Latest Handback File" (F) with links back
#   to the delivered files, and stamp "Latest Handback DateTime" (G).

$wb = $excel.ActiveWorkbook

# cornflower blue used by the workbook's existing "HyperLink" cell style
$hyperlinkColor = 15570276   # RGB(0x64, 0x95, 0xED)

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Per-language sheets: fill Latest Target File / Latest Handback File /
#    Latest Handback DateTime for the a.md and b.md rows.
# ---------------------------------------------------------------------------
$langs = @(
    @{ Sheet = "zh-cn"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf";
       HandbackCommit = "887b6899b108ac6e6006777968a13a82648f414e"; TargetCommit = "9b5047ceb2a5701f1dc843013c2bd2defd902952";
       HandbackDateTime = "2016-01-28 04:01:08" },
    @{ Sheet = "de-de"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf";
       HandbackCommit = "76d37f7a14c984c6eb78e448970f1b2a9b442c6a"; TargetCommit = "9b5047ceb2a5701f1dc843013c2bd2defd902952";
       HandbackDateTime = "2016-01-28 04:01:25" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)
    $targetMdUrl = "https://github.com/OpenLocalizationTestOrg/oltest." + $lang.Sheet + "/blob/" + $lang.TargetCommit + "/xinjiang/a.md"
    $handbackXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $lang.HandbackCommit + "/ol-handback/OpenLocalizationTestOrg/oltest." + $lang.Sheet + "/xinjiang/" + $lang.Xlf

    foreach ($row in @(2, 3)) {
        $eCell = $ws.Cells.Item($row, 5)   # column E - Latest Target File
        $fCell = $ws.Cells.Item($row, 6)   # column F - Latest Handback File
        $gCell = $ws.Cells.Item($row, 7)   # column G - Latest Handback DateTime

        $eCell.Value = "a.md"
        $ws.Hyperlinks.Add($eCell, $targetMdUrl, "", "", "a.md") | Out-Null
        Style-AsHyperlink $eCell

        $fCell.Value = $lang.Xlf
        $ws.Hyperlinks.Add($fCell, $handbackXlfUrl, "", "", $lang.Xlf) | Out-Null
        Style-AsHyperlink $fCell

        $gCell.Value = $lang.HandbackDateTime
    }
}
